$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header in B1: "qc_fail" -> "qc_fail_manual"
$ws.Range("B1").Value = "qc_fail_manual"

# Widen column B to fit the longer header text, then select B2 as the
# active cell (matches the selection seen after the edit).
$ws.Columns("B").ColumnWidth = 13.2
$ws.Range("B2").Select()
